# chore: Refactor analysis functions and update plot settings
#
# Updates the settings workbook to point at the "anca_panel_2" / Granulos
# analysis, and flips the meta_naming_scheme toggle off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# project_name: anca_panel_3 -> anca_panel_2
$ws.Range("B3").Value = "anca_panel_2"

# meta_naming_scheme: 1 (on) -> 0 (off)
$ws.Range("B13").Value = 0

# data_subsets: Monos_and_DCs -> Granulos
$ws.Range("B16").Value = "Granulos"

# Leave the selection where the edit happened, as the saved file shows.
$ws.Activate() | Out-Null
$ws.Range("B13").Select() | Out-Null
